# Apply updated inscription counts to the "Resumo Inscricoes Integrado" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E10" = 433
    "E11" = 293
    "E12" = 427
    "F12" = 228
    "H12" = 228
    "E14" = 108
    "E17" = 83
    "E23" = 177
    "F23" = 78
    "H23" = 78
    "E25" = 230
    "E26" = 129
    "E27" = 299
    "E28" = 174
    "F28" = 61
    "H28" = 61
    "E32" = 166
    "E33" = 257
    "F34" = 115
    "H34" = 115
    "E35" = 124
    "E37" = 138
    "F37" = 65
    "H37" = 65
    "E40" = 233
    "E41" = 351
    "E42" = 320
    "E44" = 274
    "E45" = 125
    "E46" = 278
    "E49" = 259
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
